$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "added multiple category option" - the Ghost of Tsushima Katana row used to
# be tagged with a single category ("Armory"); now it carries more than one,
# stored as a comma-separated list in the same Category cell.
$ws.Range("C2").Value = "Armory, Cosplay"

# Center-align the header row.
$ws.Rows("1").HorizontalAlignment = -4108  # xlCenter

# Center-align the Price column.
$ws.Columns("B").HorizontalAlignment = -4108  # xlCenter

# The Category column needs to widen a touch now that "Armory, Cosplay" is
# longer than any previous entry in that column.
$ws.Columns("C").ColumnWidth = 14.7

# Leave the selection on the cell that was just edited.
$ws.Range("C2").Select() | Out-Null
